$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the K column (column G) values for rows 2-7
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 3
